$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.503.92"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").Value = "1.796.57"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'223.03"
$ws.Range("E5").Value = "  -1.81%  "
$ws.Range("D6").Value = "'0.551"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'32.29"
$ws.Range("E8").Value = "  +2.08%  "
$ws.Range("E9").Value = "  +2.08%  "
$ws.Range("E10").Value = "  +6.39%  "
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").Value = "2.053.99"
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.803.85"
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'10.96"
$ws.Range("E14").Value = "  -4.65%  "
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "34.504.76"
$ws.Range("D17").Value = "'4.28"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").Value = "'69.07"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").Value = "'251.08"
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("D20").Value = "0.0₃0798"
$ws.Range("E20").Value = "  +7.21%  "
$ws.Range("D21").Value = "'11.07"
$ws.Range("E21").Value = "  +5.26%  "
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").Value = "'4.22"
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("D25").Value = "'161.43"
$ws.Range("E25").Value = "  +2.65%  "
$ws.Range("D26").Value = "'16.36"
$ws.Range("E26").Value = "  -1.51%  "
$ws.Range("D27").Value = "'7.13"
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").Value = "'550.36"
$ws.Range("E30").Value = "  +957.95%  "
$ws.Range("D31").Value = "'0.0524"
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("E33").Value = "  -1.07%  "
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("E35").Value = "  +1.57%  "
$ws.Range("D36").Value = "1.421.55"
$ws.Range("E36").Value = "  -2.14%  "
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0191"
$ws.Range("E38").Value = "  +1.49%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "'0.635"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("D40").Value = "'82.68"
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("D41").Value = "'0.952"
$ws.Range("E41").Value = "  +5.33%  "
$ws.Range("E42").Value = "  -3.97%  "
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").Value = "'2.13"
$ws.Range("E44").Value = "  +1.82%  "
$ws.Range("E45").Value = "  +2.86%  "
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("E47").Value = "  -2.53%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'12.28"
$ws.Range("E48").Value = "  +2.25%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "1.948.26"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").Value = "'105.43"
$ws.Range("E50").Value = "  +7.27%  "
$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = "  -0.09%  "
